$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the teacher's Name, Position, and Department on the bill header
$ws.Range("A3").Value = "নাম: Dr. Md. Alhaz Uddin "
$ws.Range("A4").Value = "পদবী: অধ্যাপক"
$ws.Range("F5").Value = "বিভাগ :গণিত"

# Label-wise bill: invigilation quantity = 1 (hour)
$ws.Range("G26").Value = 1

# Amount in words for the total bill
$ws.Range("A32").Value = "কথায়:দুই হাজার সাতশো টাকা মাত্র।"

# Move the active selection to B5 like in the final saved file
$ws.Range("B5").Select()

$wb.Save()
